$d = $word.ActiveDocument

# ===========================================================================
# Paragraph 1 (Heading3): "7.3 su_label(cols)` cat_col"
#   -> "7.3" + " " + VerbatimChar-styled "su_label(cols) cat_col"
# ===========================================================================

# Reduce the whole heading text down to just "7.3" (single plain run).
$p1Len = $d.Paragraphs(1).Range.Text.Length - 1   # exclude the paragraph mark
$headingWhole = $d.Range(0, $p1Len)
$headingWhole.Text = "7.3"

# Insert the separating space as its own run.
$spacePoint = $d.Range(3, 3)
$spacePoint.InsertAfter(" ")

# Build a VerbatimChar-styled copy of "su_label(cols) cat_col" by typing it
# inside an already-VerbatimChar-styled run (so it naturally inherits the
# character style with a proper <w:rStyle> reference), then cut it so it can
# be pasted into the heading.
$p2Start = $d.Paragraphs(2).Range.Start
$verbatimText = "su_label(cols) cat_col"
$scratchPoint = $d.Range($p2Start, $p2Start)
$scratchPoint.InsertAfter($verbatimText)
$scratchRange = $d.Range($p2Start, $p2Start + $verbatimText.Length)
$scratchRange.Cut()

# Paste the styled text right after "7.3 ".
$pastePoint = $d.Range(4, 4)
$pastePoint.Paste()

# ===========================================================================
# Paragraph 2 (SourceCode): insert a new VerbatimChar line + line break
# before the existing first line.
# ===========================================================================
$p2Start2 = $d.Paragraphs(2).Range.Start
$newLine = "      . post ``postname' (""Variable"") (""Summary label"") (""Cat_col"") (""Summary 1"")  (""Summary 2"")"

# Typing at the very start of the existing VerbatimChar run makes the new
# text inherit that same character style automatically.
$newLinePoint = $d.Range($p2Start2, $p2Start2)
$newLinePoint.InsertAfter($newLine)

# Add the textWrapping line break that separates the new line from the one
# that used to be first.
$breakPoint = $d.Range($p2Start2 + $newLine.Length, $p2Start2 + $newLine.Length)
$breakPoint.InsertBreak(6)
